$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the three new story descriptions below are written in this specific
# order (D11, then D4, then D9) so the shared-string table gets the new
# entries appended in the same order as the target workbook.

# Row 11: replace the "restart game" description with the new pop-up / banner copy, wrap the text
$ws.Range("D11").Value = "Display pop-up alerting user the game is over`nUpdate the banner with a message that says game is over press any key to start another. `nChange the color of the banner to red.`nReset the letters guessed and remaining guesses."
$ws.Range("D11").WrapText = $true

# Row 4: mark status as done ("x") and update the "Game start" story description
$ws.Range("A4").Value = "x"
$ws.Range("D4").Value = "Game starts when any key is pressed.`nBanner message color (amber) and message is changed to notify the user the the game started."
$ws.Range("D4").WrapText = $true

# Row 9: update the "When the user wins" description to mention the banner/top banner + color change
$ws.Range("D9").Value = "display the image of the band or singer in the left panel.`nDiplay a song from the band/singer on the top banner.`nChange the banner color to green`nPlay the song (audio)"

# Row 7: mark status as done ("x")
$ws.Range("A7").Value = "x"

# Row heights grow to fit the new multi-line text
$ws.Rows.Item(4).RowHeight = 43.2
$ws.Rows.Item(9).RowHeight = 57.6
$ws.Rows.Item(11).RowHeight = 57.6

# Update selection / scroll position to match the saved view
$ws.Range("A10").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
